$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(540472, "2025-10-24T18:30:00Z", 8, 12, "SV Werder Bremen", 28, "1. FC Union Berlin", 1, 0, "HomeWin")
    ,@(540470, "2025-10-25T13:30:00Z", 8, 19, "Eintracht Frankfurt", 20, "FC St. Pauli 1910", 2, 0, "HomeWin")
    ,@(540475, "2025-10-25T13:30:00Z", 8, 16, "FC Augsburg", 721, "RB Leipzig", 0, 6, "AwayWin")
    ,@(540476, "2025-10-25T13:30:00Z", 8, 2, "TSG 1899 Hoffenheim", 44, "1. FC Heidenheim 1846", 3, 1, "HomeWin")
    ,@(540477, "2025-10-25T13:30:00Z", 8, 7, "Hamburger SV", 11, "VfL Wolfsburg", 0, 1, "AwayWin")
    ,@(540474, "2025-10-25T13:45:00Z", 8, 18, "Borussia Mönchengladbach", 5, "FC Bayern München", 0, 3, "AwayWin")
    ,@(540471, "2025-10-25T16:30:00Z", 8, 4, "Borussia Dortmund", 1, "1. FC Köln", 1, 0, "HomeWin")
    ,@(540469, "2025-10-26T14:30:00Z", 8, 3, "Bayer 04 Leverkusen", 17, "SC Freiburg", 2, 0, "HomeWin")
    ,@(540473, "2025-10-26T16:30:00Z", 8, 10, "VfB Stuttgart", 15, "1. FSV Mainz 05", 2, 1, "HomeWin")
)

$startRow = 676
for ($r = 0; $r -lt $data.Count; $r++) {
    $rowValues = $data[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowValues[$c]
    }
}

$ws.Range("A1").Select()